$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 header values (column count of samples) updated
$ws.Range("B1").Value = 15
$ws.Range("C1").Value = 16
$ws.Range("D1").Value = 15
$ws.Range("E1").Value = 16

# Row 2 (CON) updated means
$ws.Range("B2").Value = 38.665747953482509
$ws.Range("C2").Value = 9.0623843370648842
$ws.Range("D2").Value = 14.227035812698375
$ws.Range("E2").Value = 3.0776763777496114

# Row 3 (STR) updated means
$ws.Range("B3").Value = 54.108535847815745
$ws.Range("C3").Value = 7.0653639578236493
$ws.Range("D3").Value = -7.0623268690707164
$ws.Range("E3").Value = 15.865637420340249

# Selection now only spans the updated columns (B:E) instead of the full B:AY range
$ws.Range("B1:E3").Select()
